# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) used four emoji values as status
# markers. Replace them with new markers:
#   old book-emoji 📘 -> warning sign ⚠️
#   old book-emoji 📕 -> text "-3"
#   old book-emoji 📙 -> text "+3"
#   old book-emoji 📗 -> check mark ✅
#
# "-3" and "+3" look like numbers, so Excel would normally convert them
# to numeric cells. They must stay as text, so they are entered with a
# leading apostrophe (forces a text/quote-prefixed cell) and then the
# cell style is reset to "Normal" right after, which clears the
# quote-prefix visual marker while keeping the cell's stored type as
# text - exactly like a user fixing the display after typing a
# text-like number into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Text

    if ($cur -eq "📘") {
        $cell.Value = "⚠️"
    } elseif ($cur -eq "📕") {
        $cell.Value = "'-3"
        $cell.Style = "Normal"
    } elseif ($cur -eq "📙") {
        $cell.Value = "'+3"
        $cell.Style = "Normal"
    } elseif ($cur -eq "📗") {
        $cell.Value = "✅"
    }
}
